$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ("5b" - Add Account to Database for Existing Users) ---
# The defect originally flagged here was incorrectly identified: Result flips
# from "Failure" to "Pass" and the defect note is cleared.
$ws.Range("E10").Value = "Pass"
$ws.Range("F10").Value = ""

# --- Row 11 (test #6 - Register Multiple Users In One Session) ---
# Crash still occurs; scenario text is clarified and the crash note gains a
# clarifying remark about the (otherwise) successful account creation.
$ws.Range("B11").Value = "Register Multiple Users In`nOne Sessions or after Log in"
$ws.Range("F11").Value = "Crash_01:`nProgram crash upon second`naccount creation`n<New Account is successfully created and saved>"

# --- Formatting: row 10 no longer highlighted as a defect ---
$ws.Range("A10").Interior.ColorIndex = -4142
$ws.Range("B10:D10").Interior.ColorIndex = -4142
$ws.Range("F10").Interior.ColorIndex = -4142
$ws.Range("B10:D10").WrapText = $true
$ws.Range("F10").WrapText = $true
$ws.Range("E10").Interior.ColorIndex = -4142

$ws.Rows.Item(10).RowHeight = 45

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 46.3

# --- Sheet view ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F12").Select()
